$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the Arduino paragraph body text with the new multi-run content.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Start = 0
$rng1.End = $d.Content.End
$ok1 = $rng1.Find.Execute('The programming of the robotic arm''s movement comes from the Arduino Uno microcontroller. With the Arduino IDE, we compiled and uploaded C code from the computer to the Arduino microcontroller using a USB cable. Fore wireless capability, the Arduino microcontroller uses the Xbee wireless module to communicate with the other Xbee connected to the SSC-32 microcontroller. With Arduino code, the SSC-32 servo controller will receive commands for changing servo positions. To determine which servo position to send to the SSC-32, it reads the xyz coordinates received from the Kinect, and converts them into servo positions with the line regression equation determined from correlating the ideal minimum and maximum range of both the Kinect coordinates and the servo positons of the SSC-32.')
if (-not $ok1) { throw "Could not find original Arduino paragraph text" }
$target1 = $d.Range($rng1.Start, $rng1.End)
$xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The programming of the robotic arm''s movement comes from the </w:t></w:r><w:r><w:t xml:space="preserve">Adafruit METRO </w:t></w:r><w:r><w:t>328 microcontroller, wh</w:t></w:r><w:r><w:t>ich is similar to an Arduino Uno</w:t></w:r><w:r><w:t xml:space="preserve">, but slightly more user-friendly. </w:t></w:r><w:r><w:t>Just like an Arduino microcontroller, it is programmable with an IDE such as the Arduino IDE</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">The four LEDs are placed on the edge of the PCB so that they are seen easier when the </w:t></w:r><w:r><w:t xml:space="preserve">METRO has a shield </w:t></w:r><w:r><w:t>mounted on top of it.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">For easier debugging, these </w:t></w:r><w:r><w:t xml:space="preserve">indicator </w:t></w:r><w:r><w:t xml:space="preserve">LEDs have one green power LED, two RX/TX LEDs, and a red LED connected to pin PF5. </w:t></w:r><w:r><w:t>Adafruit designed this microcontroller specifically</w:t></w:r><w:r><w:t xml:space="preserve"> to run the Atmega328 brain</w:t></w:r><w:r><w:t>, which has</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>32 KB of Flash and 2 KB of RAM, running at 16 MHz and preloaded with the Optiboot bootloader.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">The METRO has an FTDI USB-to-Serial converter for the purpose of sending and receiving data to a computer. </w:t></w:r><w:r><w:t>The logic level is at 5-</w:t></w:r><w:r><w:t>V, but could convert to 3.3</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">V logic if required. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so it now sits right before "METRO has a
#    shield" in the Arduino paragraph (this also removes it from its old
#    location automatically, since a document can only have one _GoBack).
# ---------------------------------------------------------------------------
$rngB = $d.Content
$rngB.Start = 0
$rngB.End = $d.Content.End
$okB = $rngB.Find.Execute('METRO has a shield')
if (-not $okB) { throw "Could not find bookmark anchor text" }
$bmRng = $d.Range($rngB.Start, $rngB.Start)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------------
# 3) Replace the Xbee paragraph's opening (tab + text) with the new
#    multi-run content (also splitting the tab into its own run).
# ---------------------------------------------------------------------------
$rng6 = $d.Content
$rng6.Start = 0
$rng6.End = $d.Content.End
$ok6 = $rng6.Find.Execute('The Xbee radio frequency module manages the wireless communication from the computer to the robotic arm.')
if (-not $ok6) { throw "Could not find original Xbee paragraph text" }
$tabStart = $rng6.Start - 1
$target6 = $d.Range($tabStart, $rng6.End)
$xml6 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">The Digi </w:t></w:r><w:r><w:t>International XBee 802.15.4 module</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">manages the wireless </w:t></w:r><w:r><w:t xml:space="preserve">data </w:t></w:r><w:r><w:t>communication from the computer to the robotic arm.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target6.InsertXML($xml6)

Write-Output "Done"
